$d = $word.ActiveDocument

# 1. Collapse the "Employee Name: {{Name}}" paragraph (currently split across
#    three runs with gramStart/gramEnd proofErr markers in between) into a
#    single run with the plain text "Employee Name: {{Name}}".
$d.Content.Find.Execute("Employee Name: {{Name}}", $false, $false, $false,
                         $false, $false, $true, 1, $false,
                         "Employee Name: {{Name}}", 2)

# 2. Bold the Gross Salary value cell ("NPR {{Salary}}") in the first table.
$t = $d.Tables.Item(1)
$cell = $t.Cell(3, 2)
$cell.Range.Bold = 1
